$d = $word.ActiveDocument

# Update the "Tempo total do processo" day-count figures per the commit diff.
$replacements = @(
    @("Concluído (91 dias)", "Concluído (95 dias)"),
    @("Assinatura Contrato (100 dias)", "Assinatura Contrato (104 dias)"),
    @("Assinatura Contrato (77 dias)", "Assinatura Contrato (81 dias)"),
    @("Assinatura Contrato (63 dias)", "Assinatura Contrato (67 dias)"),
    @("AGU (11 dias)", "AGU (15 dias)"),
    @("Total de dias 678", "Total de dias 698")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
